# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps on the last data
# row (row 11) of the "zh-cn" and "de-de" worksheets, reflecting the
# latest handback run times.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D11").Value = "2016-03-03 07:31:45"
$zhcn.Range("G11").Value = "2016-03-03 07:32:31"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D11").Value = "2016-03-03 07:31:56"
$dede.Range("G11").Value = "2016-03-03 07:32:50"
